# Update the dSF (column F) values for rows 3-10 with repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -2
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = 3
